$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- Extend the existing bordered row formatting down into the new row 17 ---
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Apply a date number format (keeps the existing border) to Start/End Date columns ---
$ws.Range("C3:D17").NumberFormat = "d-mmm-yy"

# --- New row: Sr Number 15 / Topic "Review" ---
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "Review"

# --- Start Date / End Date values for every training row ---
$ws.Cells.Item(3, 3).Value = "7/14/2025"
$ws.Cells.Item(3, 4).Value = "7/15/2025"

$ws.Cells.Item(4, 3).Value = "7/7/2025"
$ws.Cells.Item(4, 4).Value = "7/16/2025"

$ws.Cells.Item(5, 3).Value = "7/17/2025"
$ws.Cells.Item(5, 4).Value = "7/18/2025"

$ws.Cells.Item(6, 3).Value = "7/21/2025"
$ws.Cells.Item(6, 4).Value = "7/21/2025"

$ws.Cells.Item(7, 3).Value = "7/22/2025"
$ws.Cells.Item(7, 4).Value = "7/23/2025"

$ws.Cells.Item(8, 3).Value = "7/24/2025"
$ws.Cells.Item(8, 4).Value = "7/24/2025"

$ws.Cells.Item(9, 3).Value = "7/25/2025"
$ws.Cells.Item(9, 4).Value = "7/25/2025"

$ws.Cells.Item(10, 3).Value = "7/28/2025"
$ws.Cells.Item(10, 4).Value = "7/29/2025"

$ws.Cells.Item(11, 3).Value = "7/30/2025"
$ws.Cells.Item(11, 4).Value = "7/31/2025"

$ws.Cells.Item(12, 3).Value = "8/1/2025"
$ws.Cells.Item(12, 4).Value = "8/1/2025"

$ws.Cells.Item(13, 3).Value = "8/4/2025"
$ws.Cells.Item(13, 4).Value = "8/5/2025"

$ws.Cells.Item(14, 3).Value = "8/6/2025"
$ws.Cells.Item(14, 4).Value = "8/6/2025"

$ws.Cells.Item(15, 3).Value = "8/7/2025"
$ws.Cells.Item(15, 4).Value = "8/7/2025"

$ws.Cells.Item(16, 3).Value = "8/8/2025"
$ws.Cells.Item(16, 4).Value = "8/20/2025"

$ws.Cells.Item(17, 3).Value = "8/21/2025"
$ws.Cells.Item(17, 4).Value = "8/22/2025"

# --- Column widths (best approximation of the recalculated best-fit sizes) ---
$ws.Columns.Item(1).ColumnWidth = 8.833333333333334
$ws.Columns.Item(4).ColumnWidth = 8.166666666666666
$ws.Columns.Item(5).ColumnWidth = 20.5

# --- View / selection state ---
$ws.Application.ActiveWindow.Zoom = 175
$ws.Range("E3").Select() | Out-Null
